{"js": "// Username validation and confirmation\n// Replace the recipient's name, address, salutation, account number and\n// balance throughout the letter.\n\nconst replacements = [\n  [\"Luke Mcneil\", \"Eagan Duke\"],\n  [\"P.O. Box 689, 7026 Elementum, Av.\", \"893-4782 Nulla Rd.\"],\n  [\"Dear Miss Mcneil\", \"Dear Master Duke\"],\n  [\"0000079-000\", \"0000151-000\"],\n  [\"$ 7,300.00\", \"$ 4,438.00\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Username validation and confirmation\n# Replace the recipient's name, address, salutation, account number and\n# balance throughout the letter.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Luke Mcneil\", \"Eagan Duke\"),\n    @(\"P.O. Box 689, 7026 Elementum, Av.\", \"893-4782 Nulla Rd.\"),\n    @(\"Dear Miss Mcneil\", \"Dear Master Duke\"),\n    @(\"0000079-000\", \"0000151-000\"),\n    @(\"$ 7,300.00\", \"$ 4,438.00\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
